# Apply updated crypto price/volume data to the worksheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column D contains price text that can look like a plain number (e.g. "622.51").
# Force those specific cells to keep a Text number format so Excel does not
# silently convert the assigned string into a numeric value.
$textPriceCells = @("D5", "D6", "D9", "D10", "D11", "D12", "D13", "D14", "D19", "D20", "D21", "D22", "D23", "D24", "D25", "D26", "D27", "D30", "D31", "D32", "D35", "D36", "D37", "D38", "D39", "D40", "D42", "D43", "D45", "D46", "D47", "D48", "D49", "D51")
foreach ($cellRef in $textPriceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

# Cell value updates (Coin, Link, Price, Volume(1h))
$ws.Range("D2").Value = "71.095.50"
$ws.Range("E2").Value = "  +6.64%  "
$ws.Range("D3").Value = "3.668.17"
$ws.Range("E3").Value = "  +18.81%  "
$ws.Range("E4").Value = "  -0.10%  "
$ws.Range("D5").Value = "622.51"
$ws.Range("E5").Value = "  +8.20%  "
$ws.Range("D6").Value = "181.33"
$ws.Range("E6").Value = "  +2.65%  "
$ws.Range("D7").Value = "3.665.32"
$ws.Range("E7").Value = "  +18.72%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "0.538"
$ws.Range("E9").Value = "  +4.69%  "
$ws.Range("D10").Value = "0.163"
$ws.Range("E10").Value = "  +8.41%  "
$ws.Range("D11").Value = "6.63"
$ws.Range("E11").Value = "  +4.16%  "
$ws.Range("D12").Value = "0.500"
$ws.Range("E12").Value = "  +7.48%  "
$ws.Range("D13").Value = "40.59"
$ws.Range("E13").Value = "  +13.28%  "
$ws.Range("D14").Value = "0.0000255"
$ws.Range("E14").Value = "  +6.50%  "
$ws.Range("D15").Value = "4.272.01"
$ws.Range("E15").Value = "  +18.57%  "
$ws.Range("D16").Value = "71.080.34"
$ws.Range("E16").Value = "  +6.60%  "
$ws.Range("D17").Value = "3.661.80"
$ws.Range("E17").Value = "  +18.70%  "
$ws.Range("E18").Value = "  +1.41%  "
$ws.Range("D19").Value = "7.53"
$ws.Range("E19").Value = "  +8.36%  "
$ws.Range("D20").Value = "523.30"
$ws.Range("E20").Value = "  +9.00%  "
$ws.Range("D21").Value = "16.92"
$ws.Range("E21").Value = "  +1.14%  "
$ws.Range("D22").Value = "9.29"
$ws.Range("E22").Value = "  +20.30%  "
$ws.Range("D23").Value = "0.744"
$ws.Range("E23").Value = "  +8.42%  "
$ws.Range("D24").Value = "88.46"
$ws.Range("E24").Value = "  +6.15%  "
$ws.Range("D25").Value = "2.48"
$ws.Range("E25").Value = "  +12.07%  "
$ws.Range("D26").Value = "13.44"
$ws.Range("E26").Value = "  +6.68%  "
$ws.Range("D27").Value = "10.97"
$ws.Range("E27").Value = "  +8.81%  "
$ws.Range("E28").Value = "  -0.14%  "
$ws.Range("E29").Value = "  +11.64%  "
$ws.Range("D30").Value = "8.12"
$ws.Range("E30").Value = "  +2.36%  "
$ws.Range("D31").Value = "2.89"
$ws.Range("E31").Value = "  +11.62%  "
$ws.Range("D32").Value = "31.75"
$ws.Range("E32").Value = "  +13.84%  "
$ws.Range("E33").Value = "  +17.42%  "
$ws.Range("E34").Value = "  +3.91%  "
$ws.Range("D35").Value = "0.998"
$ws.Range("E35").Value = "  -0.06%  "
$ws.Range("D36").Value = "6.14"
$ws.Range("E36").Value = "  +10.48%  "
$ws.Range("D37").Value = "1.02"
$ws.Range("E37").Value = "  +8.52%  "
$ws.Range("D38").Value = "0.346"
$ws.Range("E38").Value = "  +11.95%  "
$ws.Range("D39").Value = "2.20"
$ws.Range("E39").Value = "  +10.68%  "
$ws.Range("D40").Value = "51.68"
$ws.Range("E40").Value = "  +5.58%  "
$ws.Range("E41").Value = "  +5.41%  "
$ws.Range("D42").Value = "45.45"
$ws.Range("E42").Value = "  -5.63%  "
$ws.Range("B43").Value = "Cosmos"
$ws.Range("C43").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D43").Value = "8.82"
$ws.Range("E43").Value = "  +6.45%  "
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "3.125.32"
$ws.Range("E44").Value = "  +12.37%  "
$ws.Range("D45").Value = "421.29"
$ws.Range("E45").Value = "  +13.82%  "
$ws.Range("D46").Value = "2.77"
$ws.Range("E46").Value = "  +4.38%  "
$ws.Range("D47").Value = "28.66"
$ws.Range("E47").Value = "  +16.41%  "
$ws.Range("D48").Value = "0.0371"
$ws.Range("E48").Value = "  +8.54%  "
$ws.Range("D49").Value = "139.01"
$ws.Range("E49").Value = "  +2.73%  "
$ws.Range("E50").Value = "  +0.00%  "
$ws.Range("D51").Value = "2.47"
$ws.Range("E51").Value = "  +11.32%  "
